$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string
$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 07:40"

# Row 28 - Israel
$ws.Range("B28").Value = 135043
$ws.Range("C28").Value = 1068
$ws.Range("D28").Value = 106294
$ws.Range("E28").Value = 27723

# Row 60 - Ghana
$ws.Range("B60").Value = 45012
$ws.Range("C60").Value = 143
$ws.Range("D60").Value = 43898
$ws.Range("E60").Value = 831

# Row 63 - Kirguistan
$ws.Range("B63").Value = 44487
$ws.Range("C63").Value = 29
$ws.Range("D63").Value = 40092
$ws.Range("E63").Value = 3335

# Row 64 - Uzbekistan
$ws.Range("B64").Value = 44107
$ws.Range("C64").Value = 214
$ws.Range("D64").Value = 41594
$ws.Range("E64").Value = 2157
$ws.Range("G64").Value = 4
$ws.Range("H64").Value = 356

# Row 75 - Australia
$ws.Range("B75").Value = 26373
$ws.Range("C75").Value = 51
$ws.Range("E75").Value = 2999
